$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 22.220.1.1
$ws.Range("C3").Value = 2631
$ws.Range("D3").Value = 92.8

# Row 4: Intel(R) Dual Band Wireless-AC 7265 - 19.51.48.1
$ws.Range("C4").Value = 431

# Row 5: Intel(R) Dual Band Wireless-AC 7260 - 17.15.0.5
$ws.Range("C5").Value = 273
$ws.Range("D5").Value = 98.2

# Row 6 and Row 7 swap adapter names, with new counts
$ws.Range("A6").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.32.1"
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 561
$ws.Range("D6").Value = 98.7

$ws.Range("A7").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.42.2"
$ws.Range("B7").Value = 38
$ws.Range("C7").Value = 2586
$ws.Range("D7").Value = 98.7

# Row 8: Totals
$ws.Range("B8").Value = 55
$ws.Range("C8").Value = 6482

# Row 16
$ws.Range("B16").Value = 29731

# Row 18
$ws.Range("B18").Value = 449371

# Row 20
$ws.Range("B20").Value = 77999

# Row 25
$ws.Range("B25").Value = 205276

# Row 26
$ws.Range("B26").Value = 40211

# Row 30
$ws.Range("B30").Value = 144782

# Row 36
$ws.Range("B36").Value = 122297
